$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Mixed model ANOVA" typo on rows 6 & 7 -> "Mix model ANOVA"
# (matches the value already used on row 8 for the same script3 group)
$ws.Range("D6").Value = "Mix model ANOVA"
$ws.Range("D7").Value = "Mix model ANOVA"

# Give column C an explicit width (sex-ratio rows now need the narrower script column)
$ws.Columns.Item(3).ColumnWidth = 7

# Add a new row to the table for the additional y18_sex_ratios analysis
# (binomial GLMM alongside the existing Pearson correlation summary)
$table = $ws.ListObjects.Item(1)
$newRow = $table.ListRows.Add()

$ws.Range("C13").Value = "script6"
$ws.Range("D13").Value = "GLMM with binomial dist"
$ws.Range("E13").Value = "SAS, GLIMMIX"
$ws.Range("F13").Value = "F = 6.43, df = 3, 19; P = 0.0034"

# Move the active selection to reflect where editing left off
[void]$ws.Range("E16").Select()
